$d = $word.ActiveDocument

# Locate the word "open" inside the sentence about the three open issues
# and swap it for "closed", matching the commit's intent (open -> closed).
$find = $d.Content
$found = $find.Find.Execute("open issues.", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target phrase 'open issues.' in the document."
}

$target = $find.Duplicate
# Narrow the range down to just the word "open" (4 characters) at the start
# of the matched phrase, leaving " issues." untouched.
$target.End = $target.Start + 4
$target.Text = "closed"

# Re-find the freshly inserted word so we have a Range scoped exactly to it.
$closed = $d.Content
$closed.Find.Execute("closed")

# Forcing a FormattedText round-trip on just this sub-range causes Word to
# materialize it as its own run, splitting the original single run into
# three runs: the text before "closed", "closed" itself, and the text after.
$closed.FormattedText = $closed.Duplicate
